# aggiornamento fino al 26/03
# Adds 5 new daily rows (234-238) to Sheet1, extending the data table
# from A1:D233 to A1:D238, matching the style of the last existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date style, border, bold, centered) of the last
# existing data row's A cell down into the new A234:A238 cells so the
# new dates render the same way as the rest of the column.
$ws.Range("A233").Copy() | Out-Null
$ws.Range("A234:A238").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @{ Row = 234; A = 44308; B = 1; C = 15; D = 242.600679281902 },
    @{ Row = 235; A = 44309; B = 1; C = 10; D = 161.7337861879347 },
    @{ Row = 236; A = 44310; B = 0; C = 8;  D = 129.3870289503477 },
    @{ Row = 237; A = 44311; B = 0; C = 6;  D = 97.0402717127608 },
    @{ Row = 238; A = 44312; B = 0; C = 5;  D = 80.86689309396733 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
}
